# "Made the model compatible with weekly data" — the env_info lookup table
# maps each environ_var_code to its aggregation method (reference_method)
# and report label. Rainfall figures ("rainfall" / "totprec") used to be
# summed (daily data accumulated over the period); with weekly data from
# CHAP they are now averaged like the other variables, so the
# reference_method for those two rows changes from "sum" to "mean".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("environ_info")

$ws.Range("B2").Value = "mean"   # rainfall
$ws.Range("B4").Value = "mean"   # totprec

# Reflect the saved file's last-used selection / scroll position.
$ws.Range("E8").Select()

$wb.Save()
